$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift all content one column to the left (delete the now-unused column A) ---
$ws.Columns.Item(1).Delete()

# --- Resize Table1 to match the shifted data range ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A6:G7"))

# --- Style the title cell (A2) before merging, so the merge only needs to carry
#     forward a single, already-correct style instead of creating extra ones ---
$title = $ws.Range("A2")
$title.Font.Bold = $true
$title.Font.Size = 14
$title.Font.Name = "Arial"
$title.Font.Color = 4473924
$title.HorizontalAlignment = -4108

# --- Merge the title row across the new table width ---
$ws.Range("A2:G2").Merge()

# --- New column widths for the re-laid-out table ---
$ws.Columns.Item(2).ColumnWidth = 13.333333333333334
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 13.333333333333334
$ws.Columns.Item(5).ColumnWidth = 17.666666666666668
$ws.Columns.Item(6).ColumnWidth = 30.166666666666668
$ws.Columns.Item(7).ColumnWidth = 19.833333333333332

# --- Print scaled down to 65% ---
$ws.PageSetup.Zoom = 65

# --- Reset the view: scroll back to the top-left and select E13 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("E13").Select()
